$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.2891508418386479
$ws.Cells.Item(2, 3).Value = 0.1067110064959991
$ws.Cells.Item(2, 4).Value = 0.08690141491899794
$ws.Cells.Item(2, 5).Value = 0.136257802795015
$ws.Cells.Item(2, 6).Value = 2.073173598615554
$ws.Cells.Item(2, 8).Value = 0.07973214163530429
$ws.Cells.Item(2, 9).Value = 1.50377406827981
$ws.Cells.Item(2, 11).Value = 0.3934059789598336
$ws.Cells.Item(2, 13).Value = 0.2581166958739658
$ws.Cells.Item(3, 2).Value = 0.2669171553449701
$ws.Cells.Item(3, 3).Value = 0.09709692343326992
$ws.Cells.Item(3, 4).Value = 0.0860362689062093
$ws.Cells.Item(3, 5).Value = 0.1251582388022214
$ws.Cells.Item(3, 6).Value = 2.010269201892342
$ws.Cells.Item(3, 8).Value = 0.07973214163530429
$ws.Cells.Item(3, 9).Value = 1.466216486311225
$ws.Cells.Item(3, 11).Value = 0.3615412303505536
$ws.Cells.Item(3, 13).Value = 0.2371571902339298
$ws.Cells.Item(4, 2).Value = 0.2535071141164167
$ws.Cells.Item(4, 3).Value = 0.09125632574736642
$ws.Cells.Item(4, 4).Value = 0.08549156246981227
$ws.Cells.Item(4, 5).Value = 0.1184174020750959
$ws.Cells.Item(4, 6).Value = 1.972309074682542
$ws.Cells.Item(4, 8).Value = 0.07973214163530429
$ws.Cells.Item(4, 9).Value = 1.443527407694191
$ws.Cells.Item(4, 11).Value = 0.3422739594771258
$ws.Cells.Item(4, 13).Value = 0.2244549096416435
$ws.Cells.Item(5, 2).Value = 0.2481029561336356
$ws.Cells.Item(5, 3).Value = 0.08889176286467659
$ws.Cells.Item(5, 4).Value = 0.08526616324944669
$ws.Cells.Item(5, 5).Value = 0.1156888395406028
$ws.Cells.Item(5, 6).Value = 1.957005840836345
$ws.Cells.Item(5, 8).Value = 0.07973214163530429
$ws.Cells.Item(5, 9).Value = 1.434374161633642
$ws.Cells.Item(5, 11).Value = 0.3344968874323797
$ws.Cells.Item(5, 13).Value = 0.2193201731428687
$ws.Cells.Item(6, 2).Value = 0.2472092509668187
$ws.Cells.Item(6, 3).Value = 0.08850006149637579
$ws.Cells.Item(6, 4).Value = 0.08522852797957725
$ws.Cells.Item(6, 5).Value = 0.1152368650859543
$ws.Cells.Item(6, 6).Value = 1.954474737007047
$ws.Cells.Item(6, 8).Value = 0.07973214163530429
$ws.Cells.Item(6, 9).Value = 1.432859853972829
$ws.Cells.Item(6, 11).Value = 0.33320999877661
$ws.Cells.Item(6, 13).Value = 0.2184700507953892
$ws.Cells.Item(7, 2).Value = 0.2534339868673214
$ws.Cells.Item(7, 3).Value = 0.09122437379566861
$ws.Cells.Item(7, 4).Value = 0.08548853657454458
$ws.Cells.Item(7, 5).Value = 0.1183805297003673
$ws.Cells.Item(7, 6).Value = 1.972102019826565
$ws.Cells.Item(7, 8).Value = 0.07973214163530429
$ws.Cells.Item(7, 9).Value = 1.443403589210035
$ws.Cells.Item(7, 11).Value = 0.3421687740369919
$ws.Cells.Item(7, 13).Value = 0.2243854931377101
$ws.Cells.Item(8, 2).Value = 0.2814344113429001
$ws.Cells.Item(8, 3).Value = 0.1033829862603568
$ws.Cells.Item(8, 4).Value = 0.08660589618234837
$ws.Cells.Item(8, 5).Value = 0.1324150529917816
$ws.Cells.Item(8, 6).Value = 2.051345845206285
$ws.Cells.Item(8, 8).Value = 0.07973214163530429
$ws.Cells.Item(8, 9).Value = 1.490746713266972
$ws.Cells.Item(8, 11).Value = 0.3823569544607608
$ws.Cells.Item(8, 13).Value = 0.2508549414455103
$ws.Cells.Item(9, 2).Value = 0.3382718547082391
$ws.Cells.Item(9, 3).Value = 0.1277314873181297
$ws.Cells.Item(9, 4).Value = 0.0886913522468511
$ws.Cells.Item(9, 5).Value = 0.160542042724714
$ws.Cells.Item(9, 6).Value = 2.212060838129958
$ws.Cells.Item(9, 8).Value = 0.07973214163530429
$ws.Cells.Item(9, 9).Value = 1.586568621275887
$ws.Cells.Item(9, 11).Value = 0.4635513476988251
$ws.Cells.Item(9, 13).Value = 0.3041085383462558
$ws.Cells.Item(10, 2).Value = 0.3812282067261208
$ws.Cells.Item(10, 3).Value = 0.1459438294197923
$ws.Cells.Item(10, 4).Value = 0.0901611916781917
$ws.Cells.Item(10, 5).Value = 0.1815995221151852
$ws.Cells.Item(10, 6).Value = 2.333468205380058
$ws.Cells.Item(10, 8).Value = 0.07973214163530429
$ws.Cells.Item(10, 9).Value = 1.658843116397676
$ws.Cells.Item(10, 11).Value = 0.5246972720039764
$ws.Cells.Item(10, 13).Value = 0.3440906120354725
$ws.Cells.Item(11, 2).Value = 0.401035444265716
$ws.Cells.Item(11, 3).Value = 0.154302792648167
$ws.Cells.Item(11, 4).Value = 0.09081677006498978
$ws.Cells.Item(11, 5).Value = 0.1912696009673951
$ws.Cells.Item(11, 6).Value = 2.389442485063654
$ws.Cells.Item(11, 8).Value = 0.07973214163530429
$ws.Cells.Item(11, 9).Value = 1.692142178706533
$ws.Cells.Item(11, 11).Value = 0.5528469701608287
$ws.Cells.Item(11, 13).Value = 0.3624734781092585
$ws.Cells.Item(12, 2).Value = 0.4085745302641612
$ws.Cells.Item(12, 3).Value = 0.1574790098159156
$ws.Cells.Item(12, 4).Value = 0.09106317665649044
$ws.Cells.Item(12, 5).Value = 0.1949448676130174
$ws.Cells.Item(12, 6).Value = 2.410747017387365
$ws.Cells.Item(12, 8).Value = 0.07973214163530429
$ws.Cells.Item(12, 9).Value = 1.704813066985452
$ws.Cells.Item(12, 11).Value = 0.5635551530297391
$ws.Cells.Item(12, 13).Value = 0.3694631709900307
$ws.Cells.Item(13, 2).Value = 0.4069491363544842
$ws.Cells.Item(13, 3).Value = 0.1567944683265807
$ws.Cells.Item(13, 4).Value = 0.09101019022643442
$ws.Cells.Item(13, 5).Value = 0.1941527309334461
$ws.Cells.Item(13, 6).Value = 2.406153871573139
$ws.Cells.Item(13, 8).Value = 0.07973214163530429
$ws.Cells.Item(13, 9).Value = 1.702081428713868
$ws.Cells.Item(13, 11).Value = 0.5612467896601174
$ws.Cells.Item(13, 13).Value = 0.3679565401984561
$ws.Cells.Item(14, 2).Value = 0.4016549161152057
$ws.Cells.Item(14, 3).Value = 0.1545638831150029
$ws.Cells.Item(14, 4).Value = 0.09083707891667814
$ws.Cells.Item(14, 5).Value = 0.1915716964803735
$ws.Cells.Item(14, 6).Value = 2.39119304523615
$ws.Cells.Item(14, 8).Value = 0.07973214163530429
$ws.Cells.Item(14, 9).Value = 1.69318338855436
$ws.Cells.Item(14, 11).Value = 0.5537269640737748
$ws.Cells.Item(14, 13).Value = 0.3630479503521116
$ws.Cells.Item(15, 2).Value = 0.3984170763592942
$ws.Cells.Item(15, 3).Value = 0.1531990056478207
$ws.Cells.Item(15, 4).Value = 0.09073080366626129
$ws.Cells.Item(15, 5).Value = 0.1899924959593946
$ws.Cells.Item(15, 6).Value = 2.382043241879217
$ws.Cells.Item(15, 8).Value = 0.07973214163530429
$ws.Cells.Item(15, 9).Value = 1.687741080586846
$ws.Cells.Item(15, 11).Value = 0.5491271840844547
$ws.Cells.Item(15, 13).Value = 0.3600450234657231
$ws.Cells.Item(16, 2).Value = 0.3799391336105487
$ws.Cells.Item(16, 3).Value = 0.1453990597681809
$ws.Cells.Item(16, 4).Value = 0.0901180879739556
$ws.Cells.Item(16, 5).Value = 0.1809694185856472
$ws.Cells.Item(16, 6).Value = 2.329825248187376
$ws.Cells.Item(16, 8).Value = 0.07973214163530429
$ws.Cells.Item(16, 9).Value = 1.656675477386585
$ws.Cells.Item(16, 11).Value = 0.5228643851459651
$ws.Cells.Item(16, 13).Value = 0.3428932101401116
$ws.Cells.Item(17, 2).Value = 0.3686718564832461
$ws.Cells.Item(17, 3).Value = 0.1406331649763786
$ws.Cells.Item(17, 4).Value = 0.0897388842294049
$ws.Cells.Item(17, 5).Value = 0.175457583791264
$ws.Cells.Item(17, 6).Value = 2.29798290400521
$ws.Cells.Item(17, 8).Value = 0.07973214163530429
$ws.Cells.Item(17, 9).Value = 1.637726077835751
$ws.Cells.Item(17, 11).Value = 0.5068389054073634
$ws.Cells.Item(17, 13).Value = 0.3324213697193628
$ws.Cells.Item(18, 2).Value = 0.3622162558480682
$ws.Cells.Item(18, 3).Value = 0.1378989083182489
$ws.Cells.Item(18, 4).Value = 0.08951954523982408
$ws.Cells.Item(18, 5).Value = 0.1722958598636097
$ws.Cells.Item(18, 6).Value = 2.279738130220437
$ws.Cells.Item(18, 8).Value = 0.07973214163530429
$ws.Cells.Item(18, 9).Value = 1.626866471111555
$ws.Cells.Item(18, 11).Value = 0.4976528933653981
$ws.Cells.Item(18, 13).Value = 0.3264165520435043
$ws.Cells.Item(19, 2).Value = 0.3600347937289996
$ws.Cells.Item(19, 3).Value = 0.136974325116995
$ws.Cells.Item(19, 4).Value = 0.08944506844722611
$ws.Cells.Item(19, 5).Value = 0.171226810441226
$ws.Cells.Item(19, 6).Value = 2.273572770005018
$ws.Cells.Item(19, 8).Value = 0.07973214163530429
$ws.Cells.Item(19, 9).Value = 1.623196372285861
$ws.Cells.Item(19, 11).Value = 0.4945480501413897
$ws.Cells.Item(19, 13).Value = 0.3243865536811299
$ws.Cells.Item(20, 2).Value = 0.3698686842028565
$ws.Cells.Item(20, 3).Value = 0.1411397807902688
$ws.Cells.Item(20, 4).Value = 0.08977937829623528
$ws.Cells.Item(20, 5).Value = 0.1760434419233405
$ws.Cells.Item(20, 6).Value = 2.301365312926436
$ws.Cells.Item(20, 8).Value = 0.07973214163530429
$ws.Cells.Item(20, 9).Value = 1.639739171565324
$ws.Cells.Item(20, 11).Value = 0.5085415899376926
$ws.Cells.Item(20, 13).Value = 0.3335342167026596
$ws.Cells.Item(21, 2).Value = 0.4032089097063078
$ws.Cells.Item(21, 3).Value = 0.1552187636831093
$ws.Cells.Item(21, 4).Value = 0.09088797577651775
$ws.Cells.Item(21, 5).Value = 0.1923294419643824
$ws.Cells.Item(21, 6).Value = 2.395584454954445
$ws.Cells.Item(21, 8).Value = 0.07973214163530429
$ws.Cells.Item(21, 9).Value = 1.695795290471267
$ws.Cells.Item(21, 11).Value = 0.5559343990687751
$ws.Cells.Item(21, 13).Value = 0.3644889447789268
$ws.Cells.Item(22, 2).Value = 0.4252232854566671
$ws.Cells.Item(22, 3).Value = 0.1644835867646179
$ws.Cells.Item(22, 4).Value = 0.09160175861749309
$ws.Cells.Item(22, 5).Value = 0.2030516544856411
$ws.Cells.Item(22, 6).Value = 2.457793805869954
$ws.Cells.Item(22, 8).Value = 0.07973214163530429
$ws.Cells.Item(22, 9).Value = 1.732788585509951
$ws.Cells.Item(22, 11).Value = 0.5871912726128983
$ws.Cells.Item(22, 13).Value = 0.3848860016508908
$ws.Cells.Item(23, 2).Value = 0.4134531660329515
$ws.Cells.Item(23, 3).Value = 0.1595329066418287
$ws.Cells.Item(23, 4).Value = 0.09122177261755127
$ws.Cells.Item(23, 5).Value = 0.1973217228262882
$ws.Cells.Item(23, 6).Value = 2.424533341383324
$ws.Cells.Item(23, 8).Value = 0.07973214163530429
$ws.Cells.Item(23, 9).Value = 1.713011630923461
$ws.Cells.Item(23, 11).Value = 0.5704828484825271
$ws.Cells.Item(23, 13).Value = 0.3739843250613859
$ws.Cells.Item(24, 2).Value = 0.3693275295493379
$ws.Cells.Item(24, 3).Value = 0.1409107219824648
$ws.Cells.Item(24, 4).Value = 0.08976107507205455
$ws.Cells.Item(24, 5).Value = 0.1757785533652694
$ws.Cells.Item(24, 6).Value = 2.299835933386021
$ws.Cells.Item(24, 8).Value = 0.07973214163530429
$ws.Cells.Item(24, 9).Value = 1.638828943990632
$ws.Cells.Item(24, 11).Value = 0.5077717214203403
$ws.Cells.Item(24, 13).Value = 0.3330310500839531
$ws.Cells.Item(25, 2).Value = 0.3226870816312726
$ws.Cells.Item(25, 3).Value = 0.1210888805478021
$ws.Cells.Item(25, 4).Value = 0.0881382795288026
$ws.Cells.Item(25, 5).Value = 0.1528657194833443
$ws.Cells.Item(25, 6).Value = 2.168004768815905
$ws.Cells.Item(25, 8).Value = 0.07973214163530429
$ws.Cells.Item(25, 9).Value = 1.56032120531124
$ws.Cells.Item(25, 11).Value = 0.4413268219246334
$ws.Cells.Item(25, 13).Value = 0.2895541743523253
